$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Null out the orderTimeStamp for Henley (row 9, column C)
$ws.Range("C9").ClearContents()

# Null out the orderDate for Genuine (row 7, column E)
$ws.Range("E7").ClearContents()

# Irene's (row 10) orderTimeStamp becomes a real date/time value instead of
# the unparseable text string "3/3/1900  16:25:30 PM"
$ws.Range("C10").Value = 37714.236111111109

# Update the active selection to E7
$ws.Range("E7").Select()
